$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4.389210974744354
$ws.Range("C2").Value = 0.8773251943038555
$ws.Range("B3").Value = 7.258839897500109
$ws.Range("C3").Value = 2.108458088857996
$ws.Range("B4").Value = 10.77167102313758
$ws.Range("C4").Value = 3.121451545408707
$ws.Range("B5").Value = 11.22586032083313
$ws.Range("C5").Value = 3.868091039954637
$ws.Range("B6").Value = 11.38819849730005
$ws.Range("C6").Value = 4.845529877315457
$ws.Range("B7").Value = 11.4946830466471
$ws.Range("C7").Value = 5.974366460425464
$ws.Range("B8").Value = 14.41729959441555
$ws.Range("C8").Value = 7.593048414961402
$ws.Range("B9").Value = 16.74087348140246
$ws.Range("C9").Value = 8.503492123994338
$ws.Range("B10").Value = 17.02029197930959
$ws.Range("C10").Value = 9.640255696851629
$ws.Range("B11").Value = 18.99508971640189
$ws.Range("C11").Value = 10.62161390512883
$ws.Range("B12").Value = 21.70055795608577
$ws.Range("C12").Value = 11.53614924194363
$ws.Range("B13").Value = 24.37404717765145
$ws.Range("C13").Value = 12.51286078740956
$ws.Range("B14").Value = 26.45002063519137
$ws.Range("C14").Value = 13.70581766740797
$ws.Range("B15").Value = 30.23163522335052
$ws.Range("C15").Value = 14.50634543254684
$ws.Range("B16").Value = 38.18497229427768
$ws.Range("C16").Value = 15.41443161594487
$ws.Range("B17").Value = 39.54015766885514
$ws.Range("C17").Value = 16.68896813798776
$ws.Range("B18").Value = 40.8358596361718
$ws.Range("C18").Value = 17.68843916715675
$ws.Range("B19").Value = 42.24637081988054
$ws.Range("C19").Value = 18.68784268766062
$ws.Range("B20").Value = 42.62811165670107
$ws.Range("C20").Value = 19.63036643350028
$ws.Range("B21").Value = 48.33094471918476
$ws.Range("C21").Value = 20.53304776220213
$ws.Range("B22").Value = 51.97319651995694
$ws.Range("C22").Value = 21.45714589434359
$ws.Range("B23").Value = 52.76491066245535
$ws.Range("C23").Value = 22.60393225479725
$ws.Range("B24").Value = 53.31320565783093
$ws.Range("C24").Value = 23.72180289375063
$ws.Range("B25").Value = 53.50168987103292
$ws.Range("C25").Value = 24.58782968008071
$ws.Range("B26").Value = 53.58959823329852
$ws.Range("C26").Value = 25.51752197485436
$ws.Range("B27").Value = 53.97803532637133
$ws.Range("C27").Value = 27.12191281200462
$ws.Range("B28").Value = 55.38465624084917
$ws.Range("C28").Value = 27.89521893806525
$ws.Range("B29").Value = 56.16119476823857
$ws.Range("C29").Value = 28.99267208239896
$ws.Range("B30").Value = 57.90918473792295
$ws.Range("C30").Value = 29.86617111045354
$ws.Range("B31").Value = 58.16189043055736
$ws.Range("C31").Value = 31.14006315357301
$ws.Range("B32").Value = 61.33655500163309
$ws.Range("C32").Value = 32.22828765889446
$ws.Range("B33").Value = 67.1744709698694
$ws.Range("C33").Value = 33.43726687531203
$ws.Range("B34").Value = 67.31578622706076
$ws.Range("C34").Value = 34.72665316902853
$ws.Range("B35").Value = 69.29279432400855
$ws.Range("C35").Value = 35.95087125322502
$ws.Range("B36").Value = 72.3452675463525
$ws.Range("C36").Value = 36.91516176199421
$ws.Range("B37").Value = 73.4871610513374
$ws.Range("C37").Value = 38.5491426876748
$ws.Range("B38").Value = 75.26493964864538
$ws.Range("C38").Value = 39.45277333022075
$ws.Range("B39").Value = 76.33464142600366
$ws.Range("C39").Value = 40.37614437094773
$ws.Range("B40").Value = 76.52629992203804
$ws.Range("C40").Value = 41.2919995675299
$ws.Range("B41").Value = 84.34543864074656
$ws.Range("C41").Value = 42.19171928889514
$ws.Range("B42").Value = 84.61659482284063
$ws.Range("C42").Value = 43.13107913843459
$ws.Range("B43").Value = 84.76130450755072
$ws.Range("C43").Value = 44.42179324262965
$ws.Range("B44").Value = 85.93698661043175
$ws.Range("C44").Value = 45.77655170797915
$ws.Range("B45").Value = 88.49351888909203
$ws.Range("C45").Value = 46.90526204981738
$ws.Range("B46").Value = 90.13045300567263
$ws.Range("C46").Value = 48.42265906477017
$ws.Range("B47").Value = 90.27819173470444
$ws.Range("C47").Value = 49.25298883490328
$ws.Range("B48").Value = 96.11255148235978
$ws.Range("C48").Value = 50.60341337776472
$ws.Range("B49").Value = 97.45626163485761
$ws.Range("C49").Value = 51.62124488912205

# New row 50 (appended), copy formatting from row 49 for column A
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 97.62042678732615
$ws.Range("C50").Value = 52.58098192091597
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)
